$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
$ws.Cells.Item(26, 3).Value = 23
$ws.Cells.Item(26, 4).Value = "23_Efficient Z-Gates for Quantum Computing"
$ws.Cells.Item(27, 3).Value = 24
$ws.Cells.Item(27, 4).Value = "24_Arbitrary_Waveform_Generator_for_Quantum_Informati"

$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D27").Select()
Write-Output ("ScrollRow after select: " + $excel.ActiveWindow.ScrollRow)
Write-Output ("ScrollColumn after select: " + $excel.ActiveWindow.ScrollColumn)
